$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (renamed/reordered columns) ---
$headers = @(
    "assetid", "serialnumber", "acquisitiondate", "condition", "acquisitioncost",
    "residualvalue", "usefullife", "barcode", "description", "disposaldate",
    "isconverted", "category", "location_name", "depreciationtype", "responsible_user"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Widen column A to fit the longer header/content ---
$ws.Columns.Item(1).ColumnWidth = 48

# --- Data rows ---
# columns: A assetid, B serialnumber, C acquisitiondate, D condition, E acquisitioncost,
#          F residualvalue, G usefullife, H barcode, I description, J disposaldate,
#          K isconverted, L category, M location_name, N depreciationtype, O responsible_user
$rows = @(
    @{ A = 5;  H = "AUA1000"; J = 45150.000185185185 },
    @{ A = 6;  H = "AUA1000"; J = 45150.000185185185 },
    @{ A = 9;  H = "AUA1003"; J = 45152.000185185185 },
    @{ A = 10; H = "AUA1004"; J = 45152.000185185185 },
    @{ A = 11; H = "AUA1004"; J = 45155.000185185185 },
    @{ A = 4;  H = "AUA1000"; J = 45150.000185185185 }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = "ISBN 987-897231232"
    $ws.Cells.Item($r, 3).Value2 = 45058.000185185185
    $ws.Cells.Item($r, 4).Value = "Good"
    $ws.Cells.Item($r, 5).Value = 10000
    $ws.Cells.Item($r, 6).Value = 1000
    $ws.Cells.Item($r, 7).Value = 10
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = "This would have been an epic description but nah"
    $ws.Cells.Item($r, 10).Value2 = $row.J
    $ws.Cells.Item($r, 11).Value = $false
    $ws.Cells.Item($r, 12).Value = "NewCategory"
    $ws.Cells.Item($r, 13).Value = "NewLocation4"
    $ws.Cells.Item($r, 14).Value = "Written Down Value"
    $ws.Cells.Item($r, 15).Value = "GreatestDetective"
    $r++
}

# --- Apply the date number format (numFmtId 14, "mm-dd-yy") to C2, then propagate the
#     same style to the rest of the date cells via copy/paste-special so every date
#     cell shares a single cellXfs entry (matches a single new style being added). ---
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Copy()
$ws.Range("C3:C7").PasteSpecial(-4122)
$ws.Range("J2:J7").PasteSpecial(-4122)
